$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.782.75'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.237.08'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '271.41'
$ws.Range("E5").Value = '  +5.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.62'
$ws.Range("E6").Value = '  +14.60%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +5.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.40'
$ws.Range("E10").Value = '  +8.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0961'
$ws.Range("E11").Value = '  +4.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.27'
$ws.Range("E12").Value = '  +19.11%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.570.82'
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.08'
$ws.Range("E15").Value = '  +5.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.249.57'
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.807'
$ws.Range("E17").Value = '  +4.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.780.07'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000105'
$ws.Range("E19").Value = '  +2.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.09'
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.80'
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.34'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.18'
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +7.92%  '
$ws.Range("E27").Value = '  +13.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.70'
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.54'
$ws.Range("E29").Value = '  +5.44%  '
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.70'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0923'
$ws.Range("E32").Value = '  +5.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.99'
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("E34").Value = '  +5.09%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.114'
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.124'
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0352'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.32'
$ws.Range("E38").Value = '  -3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.59'
$ws.Range("E39").Value = '  +26.68%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.228'
$ws.Range("E40").Value = '  +15.51%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.86'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '63.79'
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.36'
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0999'
$ws.Range("E45").Value = '  +1.55%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.44'
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.35'
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.446'
$ws.Range("E50").Value = '  +2.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.459.32'
$ws.Range("E51").Value = '  +2.46%  '
